$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 2 new columns before column D. This shifts the existing quarterly
# data (previously in columns D:K) two columns to the right (now F:M) and
# opens up D:E for the two newly reported quarters.
$ws.Columns("D:E").Insert()

# The freshly inserted D:E columns default to a generic style. For every row
# that holds data, copy the number formatting from the (now-shifted) F:G
# columns of that same row into D:E, so the new cells keep the same date /
# numeric formatting as the rest of the table without disturbing empty rows.
$ws.Range("F7:G7").Copy()
$ws.Range("D7:E7").PasteSpecial(-4122)
$ws.Range("F8:G8").Copy()
$ws.Range("D8:E8").PasteSpecial(-4122)
$ws.Range("F9:G9").Copy()
$ws.Range("D9:E9").PasteSpecial(-4122)
$ws.Range("F10:G10").Copy()
$ws.Range("D10:E10").PasteSpecial(-4122)
$ws.Range("F12:G12").Copy()
$ws.Range("D12:E12").PasteSpecial(-4122)
$ws.Range("F13:G13").Copy()
$ws.Range("D13:E13").PasteSpecial(-4122)
$ws.Range("F14:G14").Copy()
$ws.Range("D14:E14").PasteSpecial(-4122)
$ws.Range("F15:G15").Copy()
$ws.Range("D15:E15").PasteSpecial(-4122)
$ws.Range("F17:G17").Copy()
$ws.Range("D17:E17").PasteSpecial(-4122)
$ws.Range("F18:G18").Copy()
$ws.Range("D18:E18").PasteSpecial(-4122)
$ws.Range("F20:G20").Copy()
$ws.Range("D20:E20").PasteSpecial(-4122)
$ws.Range("F21:G21").Copy()
$ws.Range("D21:E21").PasteSpecial(-4122)
$ws.Range("F22:G22").Copy()
$ws.Range("D22:E22").PasteSpecial(-4122)
$ws.Range("F23:G23").Copy()
$ws.Range("D23:E23").PasteSpecial(-4122)
$ws.Range("F24:G24").Copy()
$ws.Range("D24:E24").PasteSpecial(-4122)
$ws.Range("F25:G25").Copy()
$ws.Range("D25:E25").PasteSpecial(-4122)
$ws.Range("F26:G26").Copy()
$ws.Range("D26:E26").PasteSpecial(-4122)
$ws.Range("F27:G27").Copy()
$ws.Range("D27:E27").PasteSpecial(-4122)
$ws.Range("F28:G28").Copy()
$ws.Range("D28:E28").PasteSpecial(-4122)
$ws.Range("F29:G29").Copy()
$ws.Range("D29:E29").PasteSpecial(-4122)
$ws.Range("F30:G30").Copy()
$ws.Range("D30:E30").PasteSpecial(-4122)
$ws.Range("F31:G31").Copy()
$ws.Range("D31:E31").PasteSpecial(-4122)
$ws.Range("F32:G32").Copy()
$ws.Range("D32:E32").PasteSpecial(-4122)
$ws.Range("F33:G33").Copy()
$ws.Range("D33:E33").PasteSpecial(-4122)
$ws.Range("F34:G34").Copy()
$ws.Range("D34:E34").PasteSpecial(-4122)
$ws.Range("F35:G35").Copy()
$ws.Range("D35:E35").PasteSpecial(-4122)
$ws.Range("F38:G38").Copy()
$ws.Range("D38:E38").PasteSpecial(-4122)
$ws.Range("F41:G41").Copy()
$ws.Range("D41:E41").PasteSpecial(-4122)
$ws.Range("F42:G42").Copy()
$ws.Range("D42:E42").PasteSpecial(-4122)
$ws.Range("F43:G43").Copy()
$ws.Range("D43:E43").PasteSpecial(-4122)
$ws.Range("F44:G44").Copy()
$ws.Range("D44:E44").PasteSpecial(-4122)
$ws.Range("F45:G45").Copy()
$ws.Range("D45:E45").PasteSpecial(-4122)
$ws.Range("F46:G46").Copy()
$ws.Range("D46:E46").PasteSpecial(-4122)
$ws.Range("F47:G47").Copy()
$ws.Range("D47:E47").PasteSpecial(-4122)
$ws.Range("F48:G48").Copy()
$ws.Range("D48:E48").PasteSpecial(-4122)
$ws.Range("F49:G49").Copy()
$ws.Range("D49:E49").PasteSpecial(-4122)
$ws.Range("F50:G50").Copy()
$ws.Range("D50:E50").PasteSpecial(-4122)
$ws.Range("F51:G51").Copy()
$ws.Range("D51:E51").PasteSpecial(-4122)
$ws.Range("F52:G52").Copy()
$ws.Range("D52:E52").PasteSpecial(-4122)
$ws.Range("F53:G53").Copy()
$ws.Range("D53:E53").PasteSpecial(-4122)
$ws.Range("F54:G54").Copy()
$ws.Range("D54:E54").PasteSpecial(-4122)
$ws.Range("F57:G57").Copy()
$ws.Range("D57:E57").PasteSpecial(-4122)
$ws.Range("F58:G58").Copy()
$ws.Range("D58:E58").PasteSpecial(-4122)
$ws.Range("F59:G59").Copy()
$ws.Range("D59:E59").PasteSpecial(-4122)
$ws.Range("F60:G60").Copy()
$ws.Range("D60:E60").PasteSpecial(-4122)
$ws.Range("F61:G61").Copy()
$ws.Range("D61:E61").PasteSpecial(-4122)
$ws.Range("F62:G62").Copy()
$ws.Range("D62:E62").PasteSpecial(-4122)
$ws.Range("F63:G63").Copy()
$ws.Range("D63:E63").PasteSpecial(-4122)
$ws.Range("F64:G64").Copy()
$ws.Range("D64:E64").PasteSpecial(-4122)
$ws.Range("F65:G65").Copy()
$ws.Range("D65:E65").PasteSpecial(-4122)
$ws.Range("F66:G66").Copy()
$ws.Range("D66:E66").PasteSpecial(-4122)
$ws.Range("F68:G68").Copy()
$ws.Range("D68:E68").PasteSpecial(-4122)
$ws.Range("F69:G69").Copy()
$ws.Range("D69:E69").PasteSpecial(-4122)
$ws.Range("F70:G70").Copy()
$ws.Range("D70:E70").PasteSpecial(-4122)
$ws.Range("F71:G71").Copy()
$ws.Range("D71:E71").PasteSpecial(-4122)
$ws.Range("F72:G72").Copy()
$ws.Range("D72:E72").PasteSpecial(-4122)
$ws.Range("F73:G73").Copy()
$ws.Range("D73:E73").PasteSpecial(-4122)
$ws.Range("F74:G74").Copy()
$ws.Range("D74:E74").PasteSpecial(-4122)
$ws.Range("F75:G75").Copy()
$ws.Range("D75:E75").PasteSpecial(-4122)
$ws.Range("F76:G76").Copy()
$ws.Range("D76:E76").PasteSpecial(-4122)
$ws.Range("F77:G77").Copy()
$ws.Range("D77:E77").PasteSpecial(-4122)
$ws.Range("F80:G80").Copy()
$ws.Range("D80:E80").PasteSpecial(-4122)
$ws.Range("F81:G81").Copy()
$ws.Range("D81:E81").PasteSpecial(-4122)
$ws.Range("F83:G83").Copy()
$ws.Range("D83:E83").PasteSpecial(-4122)
$ws.Range("F84:G84").Copy()
$ws.Range("D84:E84").PasteSpecial(-4122)
$ws.Range("F85:G85").Copy()
$ws.Range("D85:E85").PasteSpecial(-4122)
$ws.Range("F86:G86").Copy()
$ws.Range("D86:E86").PasteSpecial(-4122)
$ws.Range("F87:G87").Copy()
$ws.Range("D87:E87").PasteSpecial(-4122)
$ws.Range("F88:G88").Copy()
$ws.Range("D88:E88").PasteSpecial(-4122)
$ws.Range("F89:G89").Copy()
$ws.Range("D89:E89").PasteSpecial(-4122)
$ws.Range("F91:G91").Copy()
$ws.Range("D91:E91").PasteSpecial(-4122)
$ws.Range("F92:G92").Copy()
$ws.Range("D92:E92").PasteSpecial(-4122)
$ws.Range("F93:G93").Copy()
$ws.Range("D93:E93").PasteSpecial(-4122)
$ws.Range("F94:G94").Copy()
$ws.Range("D94:E94").PasteSpecial(-4122)
$ws.Range("F96:G96").Copy()
$ws.Range("D96:E96").PasteSpecial(-4122)
$ws.Range("F97:G97").Copy()
$ws.Range("D97:E97").PasteSpecial(-4122)
$ws.Range("F98:G98").Copy()
$ws.Range("D98:E98").PasteSpecial(-4122)
$ws.Range("F99:G99").Copy()
$ws.Range("D99:E99").PasteSpecial(-4122)
$ws.Range("F100:G100").Copy()
$ws.Range("D100:E100").PasteSpecial(-4122)
$ws.Range("F101:G101").Copy()
$ws.Range("D101:E101").PasteSpecial(-4122)
$ws.Range("F102:G102").Copy()
$ws.Range("D102:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the newly reported quarter figures (Dec-2018 in D, Sep-2018 in E)
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 133700
$ws.Range("E8").Value = 127500
$ws.Range("D9").Value = 59500
$ws.Range("E9").Value = 51300
$ws.Range("D10").Value = 74200
$ws.Range("E10").Value = 76200
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 5100
$ws.Range("E14").Value = 1500
$ws.Range("D15").Value = 26900
$ws.Range("E15").Value = 26700
$ws.Range("D17").Value = 104700
$ws.Range("E17").Value = 91300
$ws.Range("D18").Value = 29000
$ws.Range("E18").Value = 36200
$ws.Range("D20").Value = 32500
$ws.Range("E20").Value = 37200
$ws.Range("D21").Value = 88500
$ws.Range("E21").Value = 100000
$ws.Range("D22").Value = 15700
$ws.Range("E22").Value = 14900
$ws.Range("D23").Value = 45800
$ws.Range("E23").Value = 58600
$ws.Range("D24").Value = -100
$ws.Range("E24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 45800
$ws.Range("E26").Value = 58600
$ws.Range("D27").Value = 28600
$ws.Range("E27").Value = 47000
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -32500
$ws.Range("E32").Value = -37200
$ws.Range("D33").Value = 28600
$ws.Range("E33").Value = 47000
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 28600
$ws.Range("E35").Value = 47000
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 4300
$ws.Range("E41").Value = 400
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 97900
$ws.Range("E43").Value = 85500
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 4000
$ws.Range("E45").Value = 4400
$ws.Range("D46").Value = 106300
$ws.Range("E46").Value = 90200
$ws.Range("D47").Value = 649300
$ws.Range("E47").Value = 660300
$ws.Range("D48").Value = 1963700
$ws.Range("E48").Value = 1911600
$ws.Range("D49").Value = 289600
$ws.Range("E49").Value = 297400
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 11700
$ws.Range("E52").Value = 18600
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 3020600
$ws.Range("E54").Value = 2978100
$ws.Range("D57").Value = 38400
$ws.Range("E57").Value = 22600
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 71600
$ws.Range("E59").Value = 66100
$ws.Range("D60").Value = 110100
$ws.Range("E60").Value = 88700
$ws.Range("D61").Value = 1257700
$ws.Range("E61").Value = 1175300
$ws.Range("D62").Value = 431600
$ws.Range("E62").Value = 463000
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 1799300
$ws.Range("E66").Value = 1738000
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 318900
$ws.Range("E70").Value = 326100
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 0
$ws.Range("E72").Value = 0
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 902400
$ws.Range("E76").Value = 913900
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 28600
$ws.Range("E81").Value = 47000
$ws.Range("D83").Value = 27000
$ws.Range("E83").Value = 26600
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 61400
$ws.Range("E89").Value = 56400
$ws.Range("D91").Value = -63600
$ws.Range("E91").Value = -46600
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -79500
$ws.Range("E94").Value = -46500
$ws.Range("D96").Value = -45200
$ws.Range("E96").Value = -45200
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 22100
$ws.Range("E100").Value = -17700
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 4000
$ws.Range("E102").Value = -7800
